# Edit workbook to fix the typo of US import tariff (0.4 -> 0.04) on Sheet3.
$wb = $excel.ActiveWorkbook

# Sheet3 holds the values that need updating (it currently mirrors the old,
# incorrect 0.4 tariff scenario; Sheet2 already reflects the corrected one).
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Update Sheet3 cell values -------------------------------------------------
$ws3.Range("C5").Value = 0
$ws3.Range("G5").Value = 7000000
$ws3.Range("H5").Value = 5000000

$ws3.Range("C8").Value = 4000000
$ws3.Range("D8").Value = 0
$ws3.Range("E8").Value = 3000000

$ws3.Range("D12").Value = 1000000
$ws3.Range("G12").Value = 17000000
$ws3.Range("H12").Value = 0

# --- Switch the active / selected tab from Sheet2 to Sheet3 -------------------
$ws2.Select()
$ws2.Range("H13").Select()

$ws3.Select()
$ws3.Range("I3").Select()

$wb.Save()
